$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 66
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 66
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 66
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -406
$ws.Range("H39").Value = 1347.75
$ws.Range("I39").Value = 1726
$ws.Range("J39").Value = 1120.8
$ws.Range("K39").Value = 5178
$ws.Range("L39").Value = 3362.4
$ws.Range("M39").Value = -4882
$ws.Range("N39").Value = -3954.4
$ws.Range("H64").Value = 5460.778
$ws.Range("J64").Value = 5699.25
$ws.Range("L64").Value = 5699.25
$ws.Range("N64").Value = -6195.25
$ws.Range("H67").Value = 5460.778
$ws.Range("J67").Value = 5699.25
$ws.Range("L67").Value = 5699.25
$ws.Range("N67").Value = -7415.25
$ws.Range("H92").Value = 1014.63635
$ws.Range("I92").Value = 1259.1428
$ws.Range("K92").Value = 1259.1428
$ws.Range("M92").Value = -11.14280000000008
$ws.Range("H112").Value = 1631.8772
$ws.Range("J112").Value = 1665.1296
$ws.Range("L112").Value = 4995.3888
$ws.Range("N112").Value = -7211.3888
$ws.Range("H131").Value = 4377.5835
$ws.Range("I131").Value = 1718.1666
$ws.Range("J131").Value = 7037
$ws.Range("K131").Value = 5154.4998
$ws.Range("L131").Value = 21111
$ws.Range("M131").Value = -114.4997999999996
$ws.Range("N131").Value = -31191
$ws.Range("H132").Value = 1816.1587
$ws.Range("I132").Value = 1805.2034
$ws.Range("K132").Value = 5415.6102
$ws.Range("M132").Value = -2885.6102
$ws.Range("H135").Value = 916.43475
$ws.Range("I135").Value = 651.55554
$ws.Range("J135").Value = 1870
$ws.Range("K135").Value = 5863.99986
$ws.Range("L135").Value = 16830
$ws.Range("M135").Value = -3328.99986
$ws.Range("N135").Value = -21900
$ws.Range("H137").Value = 2104.3667
$ws.Range("I137").Value = 1158.878
$ws.Range("J137").Value = 4144.6313
$ws.Range("K137").Value = 3476.634
$ws.Range("L137").Value = 12433.8939
$ws.Range("M137").Value = -926.634
$ws.Range("N137").Value = -17533.8939
$ws.Range("H141").Value = 1522.8572
$ws.Range("I141").Value = 1489
$ws.Range("K141").Value = 4467
$ws.Range("M141").Value = 713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 930.7875
$ws.Range("I32").Value = 733.3521
$ws.Range("K32").Value = 733.3521
$ws.Range("M32").Value = -446.3521
$ws.Range("H74").Value = 4999
$ws.Range("I74").Value = 3844.359
$ws.Range("K74").Value = 3844.359
$ws.Range("M74").Value = -2970.359
$ws.Range("H77").Value = 4999
$ws.Range("I77").Value = 3844.359
$ws.Range("K77").Value = 19221.795
$ws.Range("M77").Value = -14853.795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 22468.27
$ws.Range("I99").Value = 26173.863
$ws.Range("J99").Value = 2087.5
$ws.Range("K99").Value = 26173.863
$ws.Range("L99").Value = 2087.5
$ws.Range("M99").Value = -24675.863
$ws.Range("N99").Value = -5083.5
$ws.Range("H105").Value = 3363.2222
$ws.Range("I105").Value = 2895.8572
$ws.Range("K105").Value = 2895.8572
$ws.Range("M105").Value = -1148.8572
$ws.Range("H107").Value = 3824.4546
$ws.Range("J107").Value = 3000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H134").Value = 2083.9285
$ws.Range("I134").Value = 1422.3513
$ws.Range("K134").Value = 4267.0539
$ws.Range("M134").Value = -1732.0539

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 209.73334
$ws.Range("I7").Value = 179.375
$ws.Range("K7").Value = 179.375
$ws.Range("M7").Value = -66.375
$ws.Range("H31").Value = 5242.353
$ws.Range("I31").Value = 10138.8
$ws.Range("J31").Value = 3202.1667
$ws.Range("K31").Value = 10138.8
$ws.Range("L31").Value = 3202.1667
$ws.Range("M31").Value = -9843.799999999999
$ws.Range("N31").Value = -3792.1667
$ws.Range("H34").Value = 5242.353
$ws.Range("I34").Value = 10138.8
$ws.Range("J34").Value = 3202.1667
$ws.Range("K34").Value = 10138.8
$ws.Range("L34").Value = 3202.1667
$ws.Range("M34").Value = -9936.799999999999
$ws.Range("N34").Value = -3606.1667
$ws.Range("H58").Value = 3821.2144
$ws.Range("I58").Value = 2033
$ws.Range("J58").Value = 5162.375
$ws.Range("K58").Value = 2033
$ws.Range("L58").Value = 5162.375
$ws.Range("M58").Value = -1830
$ws.Range("N58").Value = -5568.375
$ws.Range("H107").Value = 567.2353000000001
$ws.Range("I107").Value = 475.1905
$ws.Range("J107").Value = 715.9231
$ws.Range("K107").Value = 475.1905
$ws.Range("L107").Value = 715.9231
$ws.Range("M107").Value = 1444.8095
$ws.Range("N107").Value = -4555.9231
$ws.Range("H132").Value = 4723.6665
$ws.Range("I132").Value = 2583.1667
$ws.Range("K132").Value = 7749.500100000001
$ws.Range("M132").Value = -5219.500100000001
$ws.Range("H136").Value = 3821.2144
$ws.Range("I136").Value = 2033
$ws.Range("J136").Value = 5162.375
$ws.Range("K136").Value = 6099
$ws.Range("L136").Value = 15487.125
$ws.Range("M136").Value = -3549
$ws.Range("N136").Value = -20587.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1583
$ws.Range("J5").Value = 2851.6365
$ws.Range("L5").Value = 8554.9095
$ws.Range("N5").Value = -8778.9095
$ws.Range("H98").Value = 1064.4762
$ws.Range("I98").Value = 804.8889
$ws.Range("K98").Value = 2414.6667
$ws.Range("M98").Value = -916.6667000000002
$ws.Range("H121").Value = 818.2143
$ws.Range("I121").Value = 672
$ws.Range("K121").Value = 2016
$ws.Range("M121").Value = -706
$ws.Range("H135").Value = 1583
$ws.Range("J135").Value = 2851.6365
$ws.Range("L135").Value = 25664.7285
$ws.Range("N135").Value = -30734.7285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 18432
$ws.Range("J44").Value = 18432
$ws.Range("L44").Value = 18432
$ws.Range("N44").Value = -19624
$ws.Range("H58").Value = 23400
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 26750
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 26750
$ws.Range("N58").Value = -27304
$ws.Range("M58").Value = -9723
$ws.Range("H80").Value = 34033.863
$ws.Range("I80").Value = 46182.332
$ws.Range("K80").Value = 46182.332
$ws.Range("M80").Value = -45184.332
$ws.Range("H83").Value = 34033.863
$ws.Range("I83").Value = 46182.332
$ws.Range("K83").Value = 230911.66
$ws.Range("M83").Value = -225919.66
$ws.Range("H102").Value = 7337.0835
$ws.Range("I102").Value = 9791.923000000001
$ws.Range("J102").Value = 4435.909
$ws.Range("K102").Value = 9791.923000000001
$ws.Range("L102").Value = 4435.909
$ws.Range("M102").Value = -8169.923000000001
$ws.Range("N102").Value = -7679.909
$ws.Range("H113").Value = 5328.8423
$ws.Range("I113").Value = 5519.077
$ws.Range("J113").Value = 4916.6665
$ws.Range("K113").Value = 5519.077
$ws.Range("L113").Value = 4916.6665
$ws.Range("M113").Value = -3349.077
$ws.Range("N113").Value = -9256.666499999999
$ws.Range("H122").Value = 5960.3687
$ws.Range("J122").Value = 2328.4285
$ws.Range("L122").Value = 6985.2855
$ws.Range("N122").Value = -11885.2855
$ws.Range("H132").Value = 11691.281
$ws.Range("I132").Value = 6127.385
$ws.Range("K132").Value = 18382.155
$ws.Range("M132").Value = -15852.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9690.691999999999
$ws.Range("J16").Value = 3999.5
$ws.Range("L16").Value = 3999.5
$ws.Range("N16").Value = -4339.5
$ws.Range("H22").Value = 1143
$ws.Range("I22").Value = 1179
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 1179
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -884
$ws.Range("N22").Value = -1589
$ws.Range("H27").Value = 1143
$ws.Range("I27").Value = 1179
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 1179
$ws.Range("L27").Value = 999
$ws.Range("M27").Value = -1072
$ws.Range("N27").Value = -1213
$ws.Range("H61").Value = 2853.2856
$ws.Range("I61").Value = 2791.6956
$ws.Range("J61").Value = 3136.6
$ws.Range("K61").Value = 2791.6956
$ws.Range("L61").Value = 3136.6
$ws.Range("M61").Value = -2589.6956
$ws.Range("N61").Value = -3540.6
$ws.Range("H98").Value = 97736.664
$ws.Range("J98").Value = 97736.664
$ws.Range("L98").Value = 97736.664
$ws.Range("N98").Value = -103726.664
$ws.Range("H113").Value = 2853.2856
$ws.Range("I113").Value = 2791.6956
$ws.Range("J113").Value = 3136.6
$ws.Range("K113").Value = 2791.6956
$ws.Range("L113").Value = 3136.6
$ws.Range("M113").Value = -621.6956
$ws.Range("N113").Value = -7476.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 42239.2
$ws.Range("J63").Value = 42239.2
$ws.Range("L63").Value = 42239.2
$ws.Range("N63").Value = -43487.2
$ws.Range("H66").Value = 42239.2
$ws.Range("J66").Value = 42239.2
$ws.Range("L66").Value = 126717.6
$ws.Range("N66").Value = -132957.6
$ws.Range("H122").Value = 1565.825
$ws.Range("I122").Value = 1459.5714
$ws.Range("K122").Value = 4378.7142
$ws.Range("M122").Value = -1928.7142
$ws.Range("H132").Value = 3454.9434
$ws.Range("I132").Value = 1984.4773
$ws.Range("K132").Value = 5953.4319
$ws.Range("M132").Value = -3423.4319

